$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("BD2").Value = 126

# Row 4
$ws.Range("G4").Value = 1.91
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 5
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5

# Row 7
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6

# Row 8
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 2.03

# Row 9
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 2.7
$ws.Range("L9").Value = 3.5
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("W9").Value = 7.5
$ws.Range("AC9").Value = 8
$ws.Range("AI9").Value = 11
$ws.Range("AN9").Value = 4.5
$ws.Range("AO9").Value = 15
$ws.Range("AS9").Value = 251
$ws.Range("AV9").Value = 67
$ws.Range("AY9").Value = 29
$ws.Range("BB9").Value = 251
